$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (dSF) updates
$ws.Range("F9").Value  = -2
$ws.Range("F13").Value = 5
$ws.Range("F17").Value = 1
$ws.Range("F19").Value = -2
$ws.Range("F20").Value = -7
$ws.Range("F24").Value = -5
$ws.Range("F25").Value = -2
$ws.Range("F26").Value = -13
$ws.Range("F27").Value = -1
$ws.Range("F29").Value = 4
$ws.Range("F30").Value = 3
$ws.Range("F33").Value = 5
$ws.Range("F35").Value = -4
$ws.Range("F36").Value = -7
$ws.Range("F38").Value = 3
$ws.Range("F39").Value = -5

# Row 37 updates (E, F, H, I)
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("H37").Value = 2
$ws.Range("I37").Value = 5
